# Opencart_LoginData.xlsx edit script
# Commit: "Added more TC for register webpage"
#
# Summary of changes applied:
#  1. Sheet1 ("Sheet1"): remove the 4 mailto hyperlinks on A2:A5 and reset
#     those cells back to the plain/default look (no border, Normal font,
#     centered) that a hyperlink-style removal leaves behind; move the
#     active selection to B10.
#  2. Sheet2 is renamed to "Registration" and its sample data is replaced:
#     it used to hold bad "username" variants, it now holds bad "email"
#     format test values, with a real (clickable) hyperlink on the new
#     A3 cell; column B is widened; the active selection moves to D10.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 : drop the hyperlinks that used to sit on A2:A5
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Hyperlinks.Delete() | Out-Null

# Once the hyperlinks are gone the cells fall back to plain styling
# (default font, no border) but stay centered like the rest of the table.
$hadLinks = $ws1.Range("A2:A5")
$hadLinks.Style = "Normal"
$hadLinks.HorizontalAlignment = -4108   # xlCenter

$ws1.Activate() | Out-Null
$ws1.Range("B10").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet2 : rename to "Registration" and swap in the new email test data
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "Registration"

$ws2.Range("A1").Value = "email"
$ws2.Range("B1").Value = "res"

$ws2.Range("A2").Value = "johndoe.example.com"
$ws2.Range("B2").Value = "Invalid"

$ws2.Range("A3").Value = "johndoe@example"
$ws2.Range("B3").Value = "Invalid"

# the rest of the old sample rows are cleared out
$ws2.Range("A4").ClearContents() | Out-Null
$ws2.Range("B4").ClearContents() | Out-Null
$ws2.Range("A5:B7").ClearContents() | Out-Null

# A3 becomes a real, clickable hyperlink
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:johndoe@example.com") | Out-Null

# widen column B to fit the longer "Invalid" test values
$ws2.Columns.Item(2).ColumnWidth = 29.25

$ws2.Activate() | Out-Null
$ws2.Range("D10").Select() | Out-Null

Write-Host "Edit applied"
